$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the data rows (A2:E6) -- every cell in this
# range becomes an empty string, while row 1 (the header) is untouched.
$ws.Range("A2:E6").Value = ""
